$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 7).Value = 'servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Cells.Item(5, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(6, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(8, 7).Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(12, 7).Value = 'marwa_mustafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(13, 7).Value = 'dina.adel@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(15, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(16, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(17, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(23, 7).Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(24, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(25, 7).Value = 'youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(27, 7).Value = 'randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, hana.amr@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg'
$ws.Cells.Item(29, 7).Value = 'eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nahla91@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Cells.Item(34, 7).Value = 'servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Cells.Item(35, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg'
$ws.Cells.Item(36, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(37, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(38, 7).Value = 'AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(43, 7).Value = 'marwa_mustafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(44, 7).Value = 'dina.adel@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(46, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(48, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(54, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(55, 7).Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(56, 7).Value = 'youstina.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(58, 7).Value = 'randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, hana.amr@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg'
$ws.Cells.Item(60, 7).Value = 'eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nahla91@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Cells.Item(64, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Cells.Item(65, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(66, 7).Value = 'Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg'
$ws.Cells.Item(67, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(68, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(69, 7).Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Cells.Item(70, 7).Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Cells.Item(71, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(72, 7).Value = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(73, 7).Value = 'Arwa.elnagar@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Cells.Item(74, 7).Value = 'mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(75, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(77, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(79, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(85, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(86, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(91, 7).Value = 'eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nahla91@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Cells.Item(95, 7).Value = 'servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(97, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(98, 7).Value = 'heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(99, 7).Value = 'heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(101, 7).Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Cells.Item(102, 7).Value = 'mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(105, 7).Value = 'marwa_mustafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(106, 7).Value = 'Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Cells.Item(108, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(109, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(116, 7).Value = 'dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Cells.Item(117, 7).Value = 'dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Cells.Item(120, 7).Value = 'dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Cells.Item(121, 7).Value = 'nancy.abdelshafy@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Cells.Item(122, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(126, 7).Value = 'servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Cells.Item(127, 7).Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(128, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(129, 7).Value = 'heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(130, 7).Value = 'heba@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(131, 7).Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Cells.Item(132, 7).Value = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Cells.Item(133, 7).Value = 'mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Cells.Item(136, 7).Value = 'marwa_mustafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(137, 7).Value = 'dina.adel@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, maimustafa@med.asu.edu.eg'
$ws.Cells.Item(138, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(139, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(140, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(141, 7).Value = 'Alshimaa_khaled@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(147, 7).Value = 'dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Cells.Item(148, 7).Value = 'dalia.t.abdelmegid@med.asu.edu.eg, heba.e.mohamed@med.asu.edu.eg'
$ws.Cells.Item(151, 7).Value = 'Amr-Saeed@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
$ws.Cells.Item(152, 7).Value = 'nancy.abdelshafy@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Cells.Item(153, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(154, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Cells.Item(157, 7).Value = 'servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(158, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(159, 7).Value = 'Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg'
$ws.Cells.Item(160, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Cells.Item(161, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(162, 7).Value = 'Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg'
$ws.Cells.Item(164, 7).Value = 'mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(166, 7).Value = 'Arwa.elnagar@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Cells.Item(167, 7).Value = 'mernamahrous@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
$ws.Cells.Item(168, 7).Value = 'Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Cells.Item(170, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(172, 7).Value = 'Alshimaa_khaled@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Cells.Item(178, 7).Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(179, 7).Value = 'menna-allah.gamil@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg'
$ws.Cells.Item(182, 7).Value = 'dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Cells.Item(184, 7).Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, nahla91@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Cells.Item(185, 7).Value = 'abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg'

Write-Host "Updated $($ws.Name) G column order for 103 rows"
